# Applies the "Updated symbol list" price-refresh edit described by the diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# D-column cells hold numeric-looking strings but must stay text cells (the
# source file stores them as inline strings, not numbers). A leading apostrophe
# tells Excel to keep the entry as literal text instead of auto-converting it
# to a Number, matching the inlineStr -> inlineStr nature of the diff.
# E-column cells are plain text (coin-rank labels) and are set directly.

$ws.Range("D2").Value = "'248.08"
$ws.Range("D3").Value = "'22.48"
$ws.Range("D4").Value = "'5.232"
$ws.Range("D5").Value = "'0.05684"
$ws.Range("D6").Value = "'3.418"
$ws.Range("D7").Value = "'6.311"
$ws.Range("D8").Value = "'0.8062"
$ws.Range("D9").Value = "'0.9008"
$ws.Range("D10").Value = "'0.1424"
$ws.Range("D11").Value = "'0.07451"
$ws.Range("D13").Value = "'0.03074"
$ws.Range("D14").Value = "'0.09394"
$ws.Range("D15").Value = "'3.881"
$ws.Range("D16").Value = "'0.001590"
$ws.Range("D17").Value = "'0.04797"
$ws.Range("D18").Value = "'0.01828"
$ws.Range("D19").Value = "'0.0005802"
$ws.Range("E19").Value = "18OneONEWorstin24h"
$ws.Range("D20").Value = "'0.006409"
$ws.Range("D21").Value = "'0.005039"
$ws.Range("D22").Value = "'0.0009965"
$ws.Range("D24").Value = "'3.692"
$ws.Range("D25").Value = "'2.167"
$ws.Range("D40").Value = "'0.03955"
$ws.Range("D41").Value = "'0.006820"
$ws.Range("E41").Value = "40KickTokenKICK"
$ws.Range("D42").Value = "'0.1068"
$ws.Range("D43").Value = "'0.003200"
$ws.Range("D44").Value = "'0.008758"
$ws.Range("D45").Value = "'0.00005576"
$ws.Range("D48").Value = "'0.1397"
$ws.Range("D49").Value = "'0.00002100"
